# Routed MSS and FPGA I/O (GPIO) signals from the Mezzanine connector to SmartFusion.
#
# This rearranges the South/North Bridge GPIO rows (9-17) so that the
# SmartFusion MSS GPIO signals (SGPIO0-5) occupy C9:C15 (was previously only
# SGPIO0/SGPIO1 at C12:C13, with ALERT / I2C1 occupying C14:C16) while the
# ALERT / WAKEUP pair moves down to C16:C17, and the mezzanine/fabric side
# (F16:F17, J14:J16, M16:M17) is updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: capture the two "bold+italic" (style 18) and "bold" (style 17)
# look-and-feel blocks from the cells that currently carry them (C14/J14
# for style 17, C15/J15 for style 18) onto their new homes, BEFORE those
# source cells get overwritten later in this script.
# ---------------------------------------------------------------------

# style 18 (bold italic, blue fill) -> goes to F16, F17, M16, M17
$ws.Range("C15").Copy($ws.Range("F16"))
$ws.Range("C15").Copy($ws.Range("F17"))
$ws.Range("C15").Copy($ws.Range("M16"))
$ws.Range("C15").Copy($ws.Range("M17"))

# style 17 (bold, blue fill) -> goes to C16, J16
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("J16"))

# ---------------------------------------------------------------------
# Step 2: re-style the cells that move to plain GPIO look (style 12) or
# italic GPIO look (style 14), sourcing the style from donor cells that
# are never themselves touched by this edit.
# ---------------------------------------------------------------------

# style 12 (plain, blue fill) -> C9, C10, C14, C15
$ws.Range("F9").Copy($ws.Range("C9"))
$ws.Range("F9").Copy($ws.Range("C10"))
$ws.Range("F9").Copy($ws.Range("C14"))
$ws.Range("F9").Copy($ws.Range("C15"))

# style 14 (italic, blue fill) -> J14, J15
$ws.Range("F12").Copy($ws.Range("J14"))
$ws.Range("F12").Copy($ws.Range("J15"))

# ---------------------------------------------------------------------
# Step 3: write the final text values for every touched cell.
# ---------------------------------------------------------------------

$ws.Range("C9").Value = "SGPIO0"
$ws.Range("C10").Value = "SGPIO1"
$ws.Range("C12").Value = "SGPIO2"
$ws.Range("C13").Value = "SGPIO3"
$ws.Range("C14").Value = "SGPIO4"
$ws.Range("C15").Value = "SGPIO5"
$ws.Range("C16").Value = "ALERT"

$ws.Range("J14").Value = "UART1 TXD"
$ws.Range("J15").Value = "UART1 RXD"
$ws.Range("J16").Value = "ALERT"

$ws.Range("F16").Value = "I2C1 SCL"
$ws.Range("F17").Value = "I2C1 SDA"
$ws.Range("M16").Value = "I2C1 SCL"
$ws.Range("M17").Value = "I2C1 SDA"

# ---------------------------------------------------------------------
# Step 4: leave the same cell selected as in the authored workbook.
# ---------------------------------------------------------------------
$ws.Range("J13").Select()
